# petty-cashBook-2021.xlsx — "Update 4-Apr-2021, end of day update."
#
# The daily petty-cash ledger on "Sheet1" is rolled forward: the entries for
# 30-Mar through 2-Apr-2021 (rows 3-25) are cleared out, the opening balance
# (SALDO AWAL, E2) is updated to the new carried-forward balance, and the
# first transaction date (A3) is advanced to 4-Apr-2021 (serial 44291). The
# now-unused Keterangan (B), Debit (C) and Credit (D) entries in that block
# are removed entirely so the shared-string table drops the now-orphaned
# labels. Running-balance formulas in column E are left untouched — they
# recalculate on their own once the inputs above them are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the opening balance (carried forward from previous page) ---
$ws.Range("E2").Value = 864025

# --- Advance the first dated entry to 4-Apr-2021 ---
$ws.Range("A3").Value = 44291

# --- Clear out the now-settled Keterangan/Debit/Credit entries for the
#     rolled-over block (rows 3-25). Clear() drops the <c> element
#     entirely (value + style), matching a fully blanked-out cell. ---
$cellsToClear = "D3","B4","D4","B5","D5","B6","C6","B7","D7","B8","C8","B9","D9","B10","D10","A11","B11","D11","B12","C12","B13","D13","B14","C14","B15","C15","B16","D16","A17","B17","D17","B18","C18","B19","D19","B20","C20","B21","D21","B22","C22","B23","D23","B24","D24","A25"

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Clear()
}

# --- Restore the working view: scrolled to the top of the ledger, with
#     D4 (the next Debit entry to fill in) selected. ---
$ws.Activate() | Out-Null
$ws.Range("D4").Select() | Out-Null
